$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.037522333333333
$ws.Range("H2").Value = 9.112567
$ws.Range("I2").Value = 0.1153015356242242
$ws.Range("J2").Value = 0.1153015356242242
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05968133333333333
$ws.Range("N2").Value = 0.179044
$ws.Range("O2").Value = 0.02602747651633847
$ws.Range("P2").Value = 0.02602747651633848
$ws.Range("Q2").Value = 0.1812833828831111
$ws.Range("R2").Value = 1.631550445948
$ws.Range("S2").Value = 0.003001008010757259
$ws.Range("T2").Value = 0.00300100801075726

$ws.Range("G3").Value = 3.037522333333333
$ws.Range("H3").Value = 9.112567
$ws.Range("I3").Value = 0.1153015356242242
$ws.Range("J3").Value = 0.1153015356242242
$ws.Range("O3").Value = 0.144012433133819
$ws.Range("P3").Value = 0.144012433133819
$ws.Range("Q3").Value = 1.003057712465445
$ws.Range("R3").Value = 9.027519412188999
$ws.Range("S3").Value = 0.01660485468931024
$ws.Range("T3").Value = 0.01660485468931024

$ws.Range("G4").Value = 3.037522333333333
$ws.Range("H4").Value = 9.112567
$ws.Range("I4").Value = 0.1153015356242242
$ws.Range("J4").Value = 0.1153015356242242
$ws.Range("O4").Value = 0.8299600903498424
$ws.Range("P4").Value = 0.8299600903498425
$ws.Range("Q4").Value = 5.780736090267666
$ws.Range("R4").Value = 52.026624812409
$ws.Range("S4").Value = 0.09569567292415669
$ws.Range("T4").Value = 0.09569567292415671

$ws.Range("I5").Value = 0.325995654495798
$ws.Range("J5").Value = 0.325995654495798
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05968133333333333
$ws.Range("N5").Value = 0.179044
$ws.Range("O5").Value = 0.02602747651633847
$ws.Range("P5").Value = 0.02602747651633848
$ws.Range("Q5").Value = 0.5125482044297777
$ws.Range("R5").Value = 4.612933839867999
$ws.Range("S5").Value = 0.008484844241817774
$ws.Range("T5").Value = 0.008484844241817774

$ws.Range("I6").Value = 0.325995654495798
$ws.Range("J6").Value = 0.325995654495798
$ws.Range("O6").Value = 0.144012433133819
$ws.Range("P6").Value = 0.144012433133819
$ws.Range("S6").Value = 0.04694742739499168
$ws.Range("T6").Value = 0.04694742739499168

$ws.Range("I7").Value = 0.325995654495798
$ws.Range("J7").Value = 0.325995654495798
$ws.Range("O7").Value = 0.8299600903498424
$ws.Range("P7").Value = 0.8299600903498425
$ws.Range("S7").Value = 0.2705633828589886
$ws.Range("T7").Value = 0.2705633828589886

$ws.Range("I8").Value = 0.5587028098799778
$ws.Range("J8").Value = 0.5587028098799777
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05968133333333333
$ws.Range("N8").Value = 0.179044
$ws.Range("O8").Value = 0.02602747651633847
$ws.Range("P8").Value = 0.02602747651633848
$ws.Range("Q8").Value = 0.8784231263964445
$ws.Range("R8").Value = 7.905808137568
$ws.Range("S8").Value = 0.01454162426376344
$ws.Range("T8").Value = 0.01454162426376344

$ws.Range("I9").Value = 0.5587028098799778
$ws.Range("J9").Value = 0.5587028098799777
$ws.Range("O9").Value = 0.144012433133819
$ws.Range("P9").Value = 0.144012433133819
$ws.Range("Q9").Value = 4.860396345913778
$ws.Range("S9").Value = 0.08046015104951709
$ws.Range("T9").Value = 0.08046015104951708

$ws.Range("I10").Value = 0.5587028098799778
$ws.Range("J10").Value = 0.5587028098799777
$ws.Range("O10").Value = 0.8299600903498424
$ws.Range("P10").Value = 0.8299600903498425
$ws.Range("Q10").Value = 28.01101892808267
$ws.Range("S10").Value = 0.4637010345666972
$ws.Range("T10").Value = 0.4637010345666971
